$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45902
$ws.Range("B2").Value = 5601.17759246688
$ws.Range("C2").Value = 5077.05995295901
$ws.Range("D2").Value = 5688
$ws.Range("E2").Value = 6210.092771
$ws.Range("F2").Value = -0.0843695211614583

$ws.Range("A3").Value = 45903
$ws.Range("B3").Value = 5600.41202886578
$ws.Range("C3").Value = 5072.10543179084
$ws.Range("D3").Value = 2952
$ws.Range("E3").Value = 6209.270517
$ws.Range("F3").Value = 113.706829996877

$ws.Range("A4").Value = 45904
$ws.Range("B4").Value = 5628.9357609315
$ws.Range("C4").Value = 5063.8089826877
$ws.Range("D4").Value = 2952
$ws.Range("E4").Value = 6257.541172
$ws.Range("F4").Value = 114.183933073175

$ws.Range("A5").Value = 45905
$ws.Range("B5").Value = 5628.93576093312
$ws.Range("C5").Value = 4337.30153879838
$ws.Range("D5").Value = 2952
$ws.Range("E5").Value = 6257.541185
$ws.Range("F5").Value = 83.9127901193858

$ws.Range("A6").Value = 45906
$ws.Range("B6").Value = 1142.23741561831
$ws.Range("C6").Value = 2171.38300480519
$ws.Range("D6").Value = 2952
$ws.Range("E6").Value = 2210.626085
$ws.Range("F6").Value = 11.9904864244532

$ws.Range("A7").Value = 45907
$ws.Range("B7").Value = 970.493609902267
$ws.Range("C7").Value = 2120.11420725647
$ws.Range("D7").Value = 2952
$ws.Range("E7").Value = 1974.343581
$ws.Range("F7").Value = 7.16517409809175

$ws.Range("A8").Value = 45908
$ws.Range("B8").Value = 6110.60514061168
$ws.Range("C8").Value = 4995.23279053147
$ws.Range("D8").Value = 2952
$ws.Range("E8").Value = 6751.234904
$ws.Range("F8").Value = 111.827606413325

$ws.Range("A9").Value = 45909
$ws.Range("B9").Value = 6110.60514061168
$ws.Range("C9").Value = 5159.48888332538
$ws.Range("D9").Value = 2952
$ws.Range("E9").Value = 6751.234904
$ws.Range("F9").Value = 118.671610279738

$ws.Range("A10").Value = 45910
$ws.Range("B10").Value = 6110.60514061168
$ws.Range("C10").Value = 5137.39318932053
$ws.Range("D10").Value = 2952
$ws.Range("E10").Value = 6751.234904
$ws.Range("F10").Value = 117.750956362869

$ws.Range("A11").Value = 45911
$ws.Range("B11").Value = 6110.60514061168
$ws.Range("C11").Value = 5155.64347602902
$ws.Range("D11").Value = 2952
$ws.Range("E11").Value = 6751.234904
$ws.Range("F11").Value = 118.511384975723

$ws.Range("A12").Value = 45912
$ws.Range("B12").Value = 6110.60514061168
$ws.Range("C12").Value = 4475.57965510026
$ws.Range("D12").Value = 2952
$ws.Range("E12").Value = 6751.234904
$ws.Range("F12").Value = 90.1753924370241

$ws.Range("A13").Value = 45913
$ws.Range("B13").Value = 1154.34388801228
$ws.Range("C13").Value = 2147.68739721136
$ws.Range("D13").Value = 2952
$ws.Range("E13").Value = 2198.157881
$ws.Range("F13").Value = 9.97922459162817

$ws.Range("A14").Value = 45914
$ws.Range("B14").Value = 1022.76404167669
$ws.Range("C14").Value = 2095.48929513746
$ws.Range("D14").Value = 2952
$ws.Range("E14").Value = 2053.055031
$ws.Range("F14").Value = 7.24084518586555

$ws.Range("A15").Value = 45915
$ws.Range("B15").Value = 6078.28452290096
$ws.Range("C15").Value = 4993.87169857382
$ws.Range("D15").Value = 2952
$ws.Range("E15").Value = 6762.305097
$ws.Range("F15").Value = 113.578844694703

